$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Value updates
$ws.Range("D2").Value = 22494
$ws.Range("E2").Value = 3239
$ws.Range("F2").Value = 3239
$ws.Range("G2").Value = 3172
$ws.Range("H2").Value = 2438
$ws.Range("I2").Value = 2297
$ws.Range("J2").Value = 141
$ws.Range("K2").Value = 410097
$ws.Range("L2").Value = 377982
$ws.Range("M2").Value = 32116
$ws.Range("N2").Value = 29220
$ws.Range("O2").Value = 2896
$ws.Range("P2").Value = 6703
$ws.Range("Q2").Value = -966
$ws.Range("R2").Value = 228
$ws.Range("S2").Value = 637
$ws.Range("T2").Value = 471
$ws.Range("V2").Value = 49459
$ws.Range("W2").Value = 14.4
$ws.Range("X2").Value = 10.84
$ws.Range("Y2").Value = 8.16
$ws.Range("Z2").Value = 0.62
$ws.Range("AA2").Value = 1176.94
$ws.Range("AB2").Value = 379.15
$ws.Range("AC2").Value = 1614
$ws.Range("AD2").Value = 7
$ws.Range("AE2").Value = 20524
$ws.Range("AF2").Value = 0.55
$ws.Range("AG2").Value = 301
$ws.Range("AH2").Value = 2.67
$ws.Range("AI2").Value = 18.67
$ws.Range("AJ2").Value = 142367237
$ws.Range("D3").Value = 35077
$ws.Range("E3").Value = 3430
$ws.Range("F3").Value = 3430
$ws.Range("G3").Value = 3434
$ws.Range("H3").Value = 3083
$ws.Range("I3").Value = 2941
$ws.Range("J3").Value = 142
$ws.Range("K3").Value = 511444
$ws.Range("L3").Value = 473673
$ws.Range("M3").Value = 37772
$ws.Range("N3").Value = 34866
$ws.Range("O3").Value = 2905
$ws.Range("P3").Value = 8453
$ws.Range("Q3").Value = 1006
$ws.Range("R3").Value = -9156
$ws.Range("S3").Value = 11031
$ws.Range("T3").Value = 655
$ws.Range("V3").Value = 55357
$ws.Range("W3").Value = 9.779999999999999
$ws.Range("X3").Value = 8.789999999999999
$ws.Range("Y3").Value = 9.18
$ws.Range("Z3").Value = 0.67
$ws.Range("AA3").Value = 1254.05
$ws.Range("AB3").Value = 346.86
$ws.Range("AC3").Value = 1770
$ws.Range("AD3").Value = 5.68
$ws.Range("AE3").Value = 20624
$ws.Range("AF3").Value = 0.49
$ws.Range("AG3").Value = 280
$ws.Range("AH3").Value = 2.79
$ws.Range("AI3").Value = 14.39
$ws.Range("AJ3").Value = 169053154
$ws.Range("D4").Value = 36550
$ws.Range("E4").Value = 3869
$ws.Range("F4").Value = 3869
$ws.Range("G4").Value = 3878
$ws.Range("H4").Value = 3019
$ws.Range("I4").Value = 2877
$ws.Range("J4").Value = 142
$ws.Range("K4").Value = 534624
$ws.Range("L4").Value = 494901
$ws.Range("M4").Value = 39723
$ws.Range("N4").Value = 36807
$ws.Range("O4").Value = 2916
$ws.Range("P4").Value = 8453
$ws.Range("Q4").Value = 326
$ws.Range("R4").Value = -5454
$ws.Range("S4").Value = 5253
$ws.Range("T4").Value = 926
$ws.Range("V4").Value = 64596
$ws.Range("W4").Value = 10.59
$ws.Range("X4").Value = 8.26
$ws.Range("Y4").Value = 8.029999999999999
$ws.Range("Z4").Value = 0.58
$ws.Range("AA4").Value = 1245.87
$ws.Range("AB4").Value = 369.95
$ws.Range("AC4").Value = 1702
$ws.Range("AD4").Value = 5.74
$ws.Range("AE4").Value = 21773
$ws.Range("AF4").Value = 0.45
$ws.Range("AG4").Value = 300
$ws.Range("AH4").Value = 3.07
$ws.Range("AI4").Value = 17.63
$ws.Range("AJ4").Value = 169053154
$ws.Range("D5").Value = 36418
$ws.Range("E5").Value = 4092
$ws.Range("F5").Value = 4092
$ws.Range("G5").Value = 4103
$ws.Range("H5").Value = 3163
$ws.Range("I5").Value = 3022
$ws.Range("J5").Value = 141
$ws.Range("K5").Value = 567338
$ws.Range("L5").Value = 525380
$ws.Range("M5").Value = 41958
$ws.Range("N5").Value = 39052
$ws.Range("O5").Value = 2906
$ws.Range("P5").Value = 8457
$ws.Range("Q5").Value = 16
$ws.Range("R5").Value = -5469
$ws.Range("S5").Value = 5866
$ws.Range("T5").Value = 831
$ws.Range("V5").Value = 68012
$ws.Range("W5").Value = 11.24
$ws.Range("X5").Value = 8.69
$ws.Range("Y5").Value = 7.97
$ws.Range("Z5").Value = 0.57
$ws.Range("AA5").Value = 1252.16
$ws.Range("AB5").Value = 396.12
$ws.Range("AC5").Value = 1787
$ws.Range("AD5").Value = 5.9
$ws.Range("AE5").Value = 23088
$ws.Range("AF5").Value = 0.46
$ws.Range("AG5").Value = 340
$ws.Range("AH5").Value = 3.22
$ws.Range("AI5").Value = 19.03
$ws.Range("AJ5").Value = 169145833
$ws.Range("D6").Value = 41607
$ws.Range("E6").Value = 3366
$ws.Range("F6").Value = 3366
$ws.Range("G6").Value = 4881
$ws.Range("H6").Value = 4060
$ws.Range("I6").Value = 3835
$ws.Range("K6").Value = 649176
$ws.Range("L6").Value = 598545
$ws.Range("M6").Value = 50631
$ws.Range("N6").Value = 43444
$ws.Range("P6").Value = 8457
$ws.Range("Q6").Value = 3687
$ws.Range("R6").Value = -12117
$ws.Range("S6").Value = 6042
$ws.Range("T6").Value = 2552
$ws.Range("V6").Value = 83186
$ws.Range("W6").Value = 8.09
$ws.Range("X6").Value = 9.76
$ws.Range("Y6").Value = 9.300000000000001
$ws.Range("Z6").Value = 0.67
$ws.Range("AA6").Value = 1182.18
$ws.Range("AB6").Value = 498.66
$ws.Range("AC6").Value = 2267
$ws.Range("AD6").Value = 3.66
$ws.Range("AE6").Value = 25684
$ws.Range("AF6").Value = 0.32
$ws.Range("AI6").Value = 15.88
$ws.Range("AJ6").Value = 169145833
$ws.Range("E7").Value = 4453
$ws.Range("G7").Value = 4634
$ws.Range("H7").Value = 3545
$ws.Range("I7").Value = 3199
$ws.Range("K7").Value = 704647
$ws.Range("L7").Value = 650056
$ws.Range("M7").Value = 54591
$ws.Range("N7").Value = 46368
$ws.Range("P7").Value = 8459
$ws.Range("Y7").Value = 7.12
$ws.Range("Z7").Value = 0.52
$ws.Range("AA7").Value = 1190.77
$ws.Range("AC7").Value = 1892
$ws.Range("AD7").Value = 3.44
$ws.Range("AE7").Value = 27413
$ws.Range("AF7").Value = 0.24
$ws.Range("AG7").Value = 390
$ws.Range("AH7").Value = 6
$ws.Range("AI7").Value = 20.62
$ws.Range("E8").Value = 4633
$ws.Range("G8").Value = 4680
$ws.Range("H8").Value = 3501
$ws.Range("I8").Value = 3158
$ws.Range("K8").Value = 736137
$ws.Range("L8").Value = 678648
$ws.Range("M8").Value = 57488
$ws.Range("N8").Value = 49020
$ws.Range("P8").Value = 8459
$ws.Range("Y8").Value = 6.62
$ws.Range("Z8").Value = 0.49
$ws.Range("AA8").Value = 1180.51
$ws.Range("AC8").Value = 1867
$ws.Range("AD8").Value = 3.48
$ws.Range("AE8").Value = 28981
$ws.Range("AF8").Value = 0.22
$ws.Range("AG8").Value = 398
$ws.Range("AH8").Value = 6.13
$ws.Range("AI8").Value = 21.34
$ws.Range("E9").Value = 4794
$ws.Range("G9").Value = 4873
$ws.Range("H9").Value = 3687
$ws.Range("I9").Value = 3292
$ws.Range("K9").Value = 768807
$ws.Range("L9").Value = 708480
$ws.Range("M9").Value = 60326
$ws.Range("N9").Value = 51769
$ws.Range("P9").Value = 8459
$ws.Range("Y9").Value = 6.53
$ws.Range("Z9").Value = 0.49
$ws.Range("AA9").Value = 1174.42
$ws.Range("AC9").Value = 1946
$ws.Range("AD9").Value = 3.34
$ws.Range("AE9").Value = 30606
$ws.Range("AF9").Value = 0.21
$ws.Range("AG9").Value = 429
$ws.Range("AH9").Value = 6.6
$ws.Range("AI9").Value = 22.03

# Cell removals (columns dropped entirely for these rows)
$ws.Range("U2").ClearContents()
$ws.Range("U3").ClearContents()
$ws.Range("U4").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("U6").ClearContents()
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
